$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text number format to existing A2:A6 cells (keeps their string values)
$ws.Range("A2:A6").NumberFormat = "@"

# Add new row 7
$ws.Range("A7").Value = 1832
$ws.Range("A7").NumberFormat = "@"
$ws.Range("B7").Value = "Pesticide Manufacturing"

# Update the selection to match target state
$ws.Range("B14").Select()
